$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.743
$ws.Range("AF5").Value = 0.971
$ws.Range("AF6").Value = 0.842
$ws.Range("AF7").Value = 0.915
$ws.Range("AF8").Value = 0.883
$ws.Range("AF9").Value = 0.743
$ws.Range("AF10").Value = 0.971
$ws.Range("AF11").Value = 0.971
$ws.Range("AF12").Value = 1.265
$ws.Range("AF13").Value = 1.6
